# "newly added iAuthor TC's"
# Regenerates the per-candidate credentials (Client Id / User Name /
# Exam Password / First Name / Last Name) and Candidate ID for the
# existing rows, and appends one brand-new candidate row (row 15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range('A2').Value = 'SYCeF769'
$ws.Range('B2').Value = 231102258
$ws.Range('C2').Value = 'qxccgww29'
$ws.Range('D2').Value = 'j!23M#Eh'
$ws.Range('F2').Value = 'oJVTlTfu'
$ws.Range('G2').Value = 'mIfL'

# --- Row 3 ---
$ws.Range('A3').Value = 'jRLqW377'
$ws.Range('B3').Value = 231102257
$ws.Range('C3').Value = 'htquedi82'
$ws.Range('D3').Value = 'k%3zF!7U'
$ws.Range('F3').Value = 'KZJkOnKv'
$ws.Range('G3').Value = 'qPEy'

# --- Row 4 ---
$ws.Range('A4').Value = 'wPvjB865'
$ws.Range('B4').Value = 231102256
$ws.Range('C4').Value = 'hqwqkod19'
$ws.Range('D4').Value = 'eb&!R87K'
$ws.Range('F4').Value = 'zRcyXxst'
$ws.Range('G4').Value = 'rIQZ'

# --- Row 5 ---
$ws.Range('A5').Value = 'EFPwP221'
$ws.Range('B5').Value = 231102255
$ws.Range('C5').Value = 'yrqgslj50'
$ws.Range('D5').Value = 'F&y85Xv#'
$ws.Range('F5').Value = 'ynYWURdu'
$ws.Range('G5').Value = 'MTES'

# --- Row 6 ---
$ws.Range('A6').Value = 'UTsTi699'
$ws.Range('B6').Value = 231102254
$ws.Range('C6').Value = 'umanfgk41'
$ws.Range('D6').Value = 'gw9C2#M&'
$ws.Range('F6').Value = 'DyhZKWgN'
$ws.Range('G6').Value = 'QjGN'

# --- Row 7 ---
$ws.Range('A7').Value = 'casVL336'
$ws.Range('B7').Value = 231102253
$ws.Range('C7').Value = 'tcaxsrp63'
$ws.Range('D7').Value = 'r!S7Xu#9'
$ws.Range('F7').Value = 'BSOAYBKy'
$ws.Range('G7').Value = 'PySb'

# --- Row 8 ---
$ws.Range('A8').Value = 'vMlQd599'
$ws.Range('B8').Value = 231102252
$ws.Range('C8').Value = 'iksdpjl84'
$ws.Range('D8').Value = 'e7TaC4#%'
$ws.Range('F8').Value = 'coQVrLMJ'
$ws.Range('G8').Value = 'RTrb'

# --- Row 9 ---
$ws.Range('A9').Value = 'rciNa939'
$ws.Range('B9').Value = 231102251
$ws.Range('C9').Value = 'npaoyya66'
$ws.Range('D9').Value = 'Q!7#Nw4s'
$ws.Range('F9').Value = 'ZuWapiJN'
$ws.Range('G9').Value = 'JhVj'

# --- Row 10 ---
$ws.Range('A10').Value = 'MbqPy726'
$ws.Range('B10').Value = 231102250
$ws.Range('C10').Value = 'lzxakcl47'
$ws.Range('D10').Value = 'ge4Y#W$2'
$ws.Range('F10').Value = 'NlEPcDEe'
$ws.Range('G10').Value = 'OkMy'

# --- Row 11 ---
$ws.Range('A11').Value = 'gJZZA516'
$ws.Range('B11').Value = 231102249
$ws.Range('C11').Value = 'xqdpbar75'
$ws.Range('D11').Value = 'T8!#yx7C'
$ws.Range('F11').Value = 'tywLTPaQ'
$ws.Range('G11').Value = 'euOK'

# --- Row 12 ---
$ws.Range('A12').Value = 'kweqP129'
$ws.Range('B12').Value = 231102248
$ws.Range('C12').Value = 'maimfik57'
$ws.Range('D12').Value = 'mGe6!A2#'
$ws.Range('F12').Value = 'AqFzNIKI'
$ws.Range('G12').Value = 'vsWb'

# --- Row 13 ---
$ws.Range('A13').Value = 'TJEhz486'
$ws.Range('B13').Value = 231102247
$ws.Range('C13').Value = 'xcrdpwo42'
$ws.Range('D13').Value = 'M9vN&!6c'
$ws.Range('F13').Value = 'jDqrXVFr'
$ws.Range('G13').Value = 'YtXc'

# --- Row 14 ---
$ws.Range('A14').Value = 'QEDtU448'
$ws.Range('B14').Value = 231102246
$ws.Range('C14').Value = 'yjafghu82'
$ws.Range('D14').Value = 'tK&6M2$s'
$ws.Range('F14').Value = 'MUfzORrz'
$ws.Range('G14').Value = 'MqSi'

# --- Row 15 (brand-new candidate row) ---
$ws.Range('A15').Value = 'sGDpu301'
$ws.Range('B15').Value = 231102245
$ws.Range('C15').Value = 'txqeywb76'
$ws.Range('D15').Value = 'R7rhF$2&'
$ws.Range('E15').Value = 'MR'
$ws.Range('F15').Value = 'fwRmmbSj'
$ws.Range('G15').Value = 'SDaG'
$ws.Range('H15').Value = 'Candidate'

# Give the new row the same bordered look as the rest of the data rows.
for ($c = 1; $c -le 8; $c++) {
    $cell = $ws.Cells.Item(15, $c)
    $cell.Borders.LineStyle = 1
}

# Keep the selection/used-range in sync with the newly added row.
$null = $ws.Range('A1:H15').Select()
